# VET - Updated Monte Carlo simulations
# Flip several "include" sensitivity-variable flags from "N" to "Y" on the
# PowerPlants and Fuels sheets, and move the active sheet/selection from
# Globals -> PowerPlants (cell B8).

$wb = $excel.ActiveWorkbook

$powerPlants = $wb.Worksheets.Item("PowerPlants")
$fuels = $wb.Worksheets.Item("Fuels")

# PowerPlants: HeatRate, DiscountRate, MaxCapacity, MaxActivity -> included
$powerPlants.Range("B7").Value = "Y"
$powerPlants.Range("B15").Value = "Y"
$powerPlants.Range("B17").Value = "Y"
$powerPlants.Range("B18").Value = "Y"

# Fuels: MaxActivity, MaxCapacity -> included
$fuels.Range("B11").Value = "Y"
$fuels.Range("B12").Value = "Y"

# Make PowerPlants the active sheet/tab, with B8 selected (Globals loses
# the tabSelected flag as a result).
$powerPlants.Activate()
$powerPlants.Range("B8").Select()
